$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new log entry as row 17
$ws.Range("A17").Value = 20170217
$ws.Range("B17").Value = "dcstarne"
$ws.Range("C17").Value = "Preassessment"
$ws.Range("D17").Value = "gcc "
$ws.Range("E17").Value = "Execute: ./main input.txt output.txt"

# Match the original row height/formatting convention used by the other data rows
$ws.Rows.Item(17).RowHeight = 15.75

# Move the sheet's active selection from E18 to D18
$ws.Range("D18").Select()
